# Consolidate the text runs in the caption textbox on slide 1.
#
# Original runs: "The" | " " | "picture" | " " | "first"
# Target runs:   "The " | "picture " | "first"
#
# Re-assigning .Text on a Characters() sub-range spanning what used to be
# two separate runs merges them into a single run, which is exactly the
# "consolidate text run nodes" behaviour described by the commit.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)          # "TextBox 3" containing "The picture first"
$tr  = $shp.TextFrame.TextRange

# Merge "The" + " " (characters 1-4) into a single run "The ".
$tr.Characters(1, 4).Text = "The "

# Merge "picture" + " " (characters 5-12) into a single run "picture ".
$tr.Characters(5, 8).Text = "picture "

# The trailing "first" run (characters 13-17) is left untouched.
